$wb = $excel.ActiveWorkbook

# --- Sheet: analysis_details -----------------------------------------------
$ws3 = $wb.Worksheets.Item("analysis_details")

# Insert a new row 3 ("focus" / 10500) later (see below); first lay down the
# new "increment / step / configure view axes" block below the existing
# table (rows 10-15, then 17) in the same order the original author typed
# them, so the shared-string table ends up in the matching order.
$ws3.Rows("3:3").Insert()

$ws3.Range("A10").Value = "increment"
$ws3.Range("B10").Value = 100

$ws3.Range("A11").Value = "step"
$ws3.Range("B11").Value = 10

$ws3.Range("A12").Value = "configure view axes"

$ws3.Range("C13").Value = "scale"
$ws3.Range("D13").Value = "15 sec"
$ws3.Range("C14").Value = "position "
$ws3.Range("D15").Value = "1400 Hz"

$ws3.Range("B13").Value = "Time"
$ws3.Range("B14").Value = "Freq"
$ws3.Range("D14").Value = "20 Hz"

$ws3.Range("C15").Value = "scale"

$ws3.Range("A3").Value = "focus"
$ws3.Range("B3").Value = 10500

$ws3.Range("A17").Value = "All saved as Use_fish_calls preset"

# Update the remembered selection for this sheet.
$ws3.Range("D20").Select()

# --- Sheet: files_to_evaluate_quiet_061920 ----------------------------------
$ws1 = $wb.Worksheets.Item("files_to_evaluate_quiet_061920")
$ws1.Range("I10").Select()
$ws1.PageSetup.Orientation = 1

# --- Restore the originally active sheet/tab (start_times) -----------------
$ws2 = $wb.Worksheets.Item("start_times")
$ws2.Activate()
